$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values: B2 and D2 updated, C2 and E2 cleared
$ws.Range("B2").Value = 5.6988614982962673
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.552842783658944
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 4.8610006685951728
$ws.Range("C3").Value = 6.057254819772349
$ws.Range("D3").Value = 3.585621399859106
$ws.Range("E3").Value = 9.1862437576153866

# Update selection to match new used-range highlight
$ws.Range("B1:E3").Select()
